$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row text changes
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Từ vựng / Cụm từ (English)"
$ws.Range("B1").Value = "Nghĩa / Giải thích (Vietnamese)"
$ws.Range("C1").Value = "Examples (Ví dụ)"

# ---------------------------------------------------------------------
# 2. B12 text change (drop the parenthetical note)
# ---------------------------------------------------------------------
$ws.Range("B12").Value = "sách tự lực"

# ---------------------------------------------------------------------
# 3. Column C: example sentences with a bolded key term (rich text)
# ---------------------------------------------------------------------
function Set-Example($addr, $prefix, $bold, $suffix) {
    # Note: the prefix run is deliberately left with NO explicit font
    # properties (matches what real Excel produces: the untouched leading
    # run inherits the cell/default font and carries no <rPr/>). Only the
    # bolded keyword and the trailing suffix get explicit rPr.
    $cell = $ws.Range($addr)
    $full = $prefix + $bold + $suffix
    $cell.Value = $full
    $boldStart = $prefix.Length + 1
    $boldLen = $bold.Length
    $b = $cell.Characters($boldStart, $boldLen)
    $b.Font.Name = "Arial"
    $b.Font.Size = 11
    $b.Font.Bold = $true
    $b.Font.Color = 2039583
    if ($suffix.Length -gt 0) {
        $post = $cell.Characters($boldStart + $boldLen, $suffix.Length)
        $post.Font.Name = "Arial"
        $post.Font.Size = 11
        $post.Font.Bold = $false
        $post.Font.Color = 2039583
    }
}

Set-Example "C2"  '"Harry Potter" is a famous fantasy ' "novel" "."
Set-Example "C3"  "She loves reading a good " "mystery" " before bed."
Set-Example "C4"  "The movie was a psychological " "thriller" " that kept us guessing."
Set-Example "C5"  "My grandmother enjoys reading a " "romance novel" " on the weekend."
Set-Example "C6"  '"Dune" is a classic ' "science fiction book" "."
Set-Example "C7"  "He published a collection of " "short stories" "."
Set-Example "C8"  "I am reading a " "biography" " of Albert Einstein."
Set-Example "C9"  "In his " "autobiography" ", he describes his childhood struggles."
Set-Example "C10" "This " "travel book" " lists the best hotels in Paris."
Set-Example "C11" "The former president wrote his " "memoirs" " after leaving office."
Set-Example "C12" "She bought a " "self-help" " book to improve her confidence."

# ---------------------------------------------------------------------
# 4. Styling: build one reference cell per target look, then propagate
#    with Copy + PasteSpecial(xlPasteFormats) to avoid style-table churn.
# ---------------------------------------------------------------------

# -- Reference cell for the "bold" style (column A all rows + B1/C1) --
$boldRef = $ws.Range("A1")
$boldRef.Font.Name = "Arial"
$boldRef.Font.Size = 11
$boldRef.Font.Bold = $true
$boldRef.Font.Color = 2039583
$boldRef.Borders.LineStyle = 1
$boldRef.Borders.Weight = -4138
$boldRef.Borders.Color = 0
$boldRef.HorizontalAlignment = -4131
$boldRef.VerticalAlignment = -4108
$boldRef.WrapText = $true
$boldRef.IndentLevel = 1
$boldRef.ReadingOrder = 1

$boldRef.Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)
$boldRef.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)

# -- Reference cell for the "normal" style (columns B & C, rows 2-12) --
$normRef = $ws.Range("B2")
$normRef.Font.Name = "Arial"
$normRef.Font.Size = 11
$normRef.Font.Bold = $false
$normRef.Font.Color = 2039583
$normRef.Borders.LineStyle = 1
$normRef.Borders.Weight = -4138
$normRef.Borders.Color = 0
$normRef.HorizontalAlignment = -4131
$normRef.VerticalAlignment = -4108
$normRef.WrapText = $true
$normRef.IndentLevel = 1
$normRef.ReadingOrder = 1

$normRef.Copy()
$ws.Range("B2:C12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(2).ColumnWidth = 34
$ws.Columns.Item(3).ColumnWidth = 53

# ---------------------------------------------------------------------
# 6. Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30.75
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 15.75
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30

# ---------------------------------------------------------------------
# 7. Selection
# ---------------------------------------------------------------------
$ws.Range("G6").Select()
